$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old rows 2 and 3 (ECs sender rows), keep what were rows 4 and 5
# (FAPs sender rows) which move up to become the new rows 2 and 3, with
# updated TPM-derived values.
$ws.Rows.Item(2).Delete() | Out-Null
$ws.Rows.Item(2).Delete() | Out-Null

# Row 2 (was row 4): FAPs, Tnfsf11, Tnfrsf11a, ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Tnfsf11"
$ws.Range("C2").Value = "Tnfrsf11a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.07928633333333333
$ws.Range("H2").Value = 0.237859
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.7227396666666666
$ws.Range("N2").Value = 2.168219
$ws.Range("O2").Value = 0.1904947117138812
$ws.Range("P2").Value = 0.1904947117138812
$ws.Range("Q2").Value = 0.05730337812455555
$ws.Range("R2").Value = 0.5157304031209999
$ws.Range("S2").Value = 0.1904947117138812
$ws.Range("T2").Value = 0.1904947117138812

# Row 3 (was row 5): FAPs, Tnfsf11, Tnfrsf11a, MuSCs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Tnfsf11"
$ws.Range("C3").Value = "Tnfrsf11a"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.07928633333333333
$ws.Range("H3").Value = 0.237859
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.071274666666666
$ws.Range("N3").Value = 9.213823999999999
$ws.Range("O3").Value = 0.8095052882861188
$ws.Range("P3").Value = 0.8095052882861187
$ws.Range("Q3").Value = 0.2435101069795555
$ws.Range("R3").Value = 2.191590962816
$ws.Range("S3").Value = 0.8095052882861188
$ws.Range("T3").Value = 0.8095052882861187
